$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in F14 from 21 to 30
$ws.Range("F14").Value = 30

# Remove rows 15-17 entirely (the SUM formula row, the 300 value row,
# and the F16-F15 difference formula row)
$ws.Range("A15:A17").EntireRow.Delete() | Out-Null

# Keep the same cell selected as before the edit
$ws.Range("F15").Select() | Out-Null
